# The intro paragraph below the campaign dates was originally copied from
# a "Perseu(ide)" activity guide and still referenced the Perseids, with the
# text also fragmented across many small runs. Replace the whole paragraph's
# text with a single corrected, consolidated run naming the Gemeni
# constellation instead.
$d = $word.ActiveDocument

$old = "Prin această activitate participați în cadrul unei campanii globale de observare și consemnare a celor mai slabe stele vizibile ca metodă de măsurare a poluării luminoase dintr-un anumit loc. Localizând și observând constelația Perseu pe cerul nopții și comparând-o cu diagramele stelare, oamenii din întreaga lume vor putea afla în ce măsură iluminatul nocturn din comunitatea lor contribuie la poluarea luminoasă. Contribuțiile dumneavoastră la baza de date online vor facilita o documentare globală privind cerul nocturn observabil."
$new = "Prin această activitate participați în cadrul unei campanii globale de observare și consemnare a celor mai slabe stele vizibile ca metodă de măsurare a poluării luminoase dintr-un anumit loc. Localizând și observând  Constelația Gemeni pe cerul nopții și comparând-o cu diagramele stelare, oamenii din întreaga lume vor putea afla în ce măsură iluminatul nocturn din comunitatea lor contribuie la poluarea luminoasă. Contribuțiile dumneavoastră la baza de date online vor facilita o documentare globală privind cerul nocturn observabil."

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Delete the matched range (all the old runs) and type the replacement
    # text fresh, which produces a single consolidated run.
    $rng.Delete()
    $rng.InsertAfter($new)
} else {
    Write-Output "WARNING: target paragraph text not found"
}
